# Cyprus Division 1 - database update (11-04-2024)
# 1) Swap a handful of mis-ordered match rows (same id sequence, rows exchange
#    their B..AC contents while keeping column A's running id unchanged).
# 2) Append 7 brand-new (not-yet-played) fixtures as rows 233..239.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: swap rows that were recorded in the wrong order
# ---------------------------------------------------------------------------
$swapPairs = @(
  @(39,40),
  @(77,78),
  @(115,116),
  @(117,118),
  @(138,139),
  @(142,143),
  @(149,150),
  @(152,153),
  @(156,157),
  @(159,160),
  @(179,180)
)

foreach ($pair in $swapPairs) {
  $r1 = $pair[0]
  $r2 = $pair[1]
  $addr1 = "B{0}:AC{0}" -f $r1
  $addr2 = "B{0}:AC{0}" -f $r2
  $rng1 = $ws.Range($addr1)
  $rng2 = $ws.Range($addr2)
  $v1 = $rng1.Value2
  $v2 = $rng2.Value2
  $rng1.Value2 = $v2
  $rng2.Value2 = $v1
}

# ---------------------------------------------------------------------------
# Part 2: append the new fixtures (rows 233..239)
# ---------------------------------------------------------------------------
$newRows = @(
  @(233, 231, 7968566, "AEL Limassol",          "Othellos Athienou",     45394.54166666666, 1.615, 3.6, 5,     1.615, 3.75, 6,     -0.75, 1.825, 1.975, 2.75, 2,     1.8),
  @(234, 232, 7968567, "Nea Salamis Famagusta",   "Apollon Limassol",     45395.45833333334, 3,     3.4, 2.1,   3.4,   3.4,  2.15,  0.25,  1.95,  1.85,  2.5,  1.875, 1.925),
  @(235, 233, 7968568, "APK Karmotissa",         "Doxa Katokopias",       45395.5,           1.833, 3.4, 3.75,  1.833, 3.6,  4.333, -0.5,  1.875, 1.925, 2.5,  1.825, 1.975),
  @(236, 234, 7968569, "Ethnikos Achnas",        "AE Zakakiou",           45395.54166666666, 1.333, 4.2, 9.5,   1.4,   4.5,  9,     -1.25, 1.8,   2,     3.25, 1.975, 1.825),
  @(237, 235, 7968552, "Pafos FC",                "Apoel Nicosia",        45396.45833333334, 2.5,   3.5, 2.4,   2.7,   3.75, 2.4,   0,     2.025, 1.775, 2.75, 1.9,   1.9),
  @(238, 236, 7968264, "Aris Limassol",           "AEK Larnaca",          45396.54166666666, 2.6,   3,   2.6,   2.75,  3.1,  2.75,  0,     1.9,   1.9,   2.25, 1.9,   1.9),
  @(239, 237, 7968553, "Anorthosis Famagusta",    "Omonia Nicosia",       45396.5625,        2.9,   3.1, 2.3,   3.75,  3.5,  2,     0.5,   1.8,   2,     2.5,  1.925, 1.875)
)

foreach ($row in $newRows) {
  $r      = $row[0]
  $idVal  = $row[1]
  $bVal   = $row[2]
  $fVal   = $row[3]
  $gVal   = $row[4]
  $eVal   = $row[5]

  $ws.Cells.Item($r, 1).Value2 = $idVal          # A - id
  $ws.Cells.Item($r, 2).Value2 = $bVal            # B
  $ws.Cells.Item($r, 3).Value2 = "Cyprus Division 1"  # C
  $ws.Cells.Item($r, 4).Value2 = "Cyprus Division 1"  # D
  $ws.Cells.Item($r, 5).Value2 = $eVal            # E - date
  $ws.Cells.Item($r, 6).Value2 = $fVal            # F - HomeTeam
  $ws.Cells.Item($r, 7).Value2 = $gVal            # G - AwayTeam

  # K..V : odds (12 values)
  for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item($r, 11 + $i).Value2 = $row[6 + $i]
  }

  # W..AA : all zero (placeholders until the match is played)
  for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($r, 23 + $i).Value2 = 0
  }
}

# Apply the same formatting used by the existing data rows (column A id style,
# column E date style) to the freshly appended rows.
$ws.Range("A232").Copy()
$ws.Range("A233:A239").PasteSpecial(-4122)
$ws.Range("E232").Copy()
$ws.Range("E233:E239").PasteSpecial(-4122)

$excel.CutCopyMode = 0
